$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: "Tabellen ueberarbeiten" status changed from "offen" to "done" ---
$ws.Range("B3").Value2 = "done"
$ws.Range("B2").Copy() | Out-Null
$ws.Range("B3").PasteSpecial(-4122) | Out-Null

# --- New Row 13: "Logging-System" task (added first, so that the new shared
#     strings end up appended in the same order as the target workbook) ---
$ws.Range("A13").Value2 = "Logging-System"
$ws.Range("B13").Value2 = "in Arbeit"
$ws.Range("B10").Copy() | Out-Null
$ws.Range("B13").PasteSpecial(-4122) | Out-Null
$ws.Range("C13").Value2 = "Jonas"

# --- Row 7: clarify task description ---
$ws.Range("A7").Value2 = "Einteilung der Produkte in Kategorien -> neue Spalte in DB-Tabelle"

# --- New Row 14: "Bilder neben Produkten anzeigen" task ---
$ws.Range("A14").Value2 = "Bilder neben Produkten anzeigen -> neue Spalte in DB-Tabelle, die Pfad zum Bild anzeigt"
$ws.Range("B14").Value2 = "offen"
$ws.Range("B9").Copy() | Out-Null
$ws.Range("B14").PasteSpecial(-4122) | Out-Null
$ws.Range("C14").Value2 = "Jonas"

# --- Row 12: mark "Login zurueckgehen in Login-Page verhindern" as done ---
$ws.Range("B12").Value2 = "done"
$ws.Range("B11").Copy() | Out-Null
$ws.Range("B12").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# --- Update selection to match the final cursor position ---
$ws.Range("A14").Select() | Out-Null
